# Generate Report for handback
# Updates the "zh-cn" and "de-de" localization-status sheets:
#   - Status column (B) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" for the two real file rows.
#   - Two new columns are populated for those rows:
#       E = Latest Target File    (same file reference as column A)
#       F = Latest Handback File  (same file reference as column C)
#     both rendered/linked exactly like their source columns.
#   - de-de additionally records the real handback timestamp in column G
#     (Latest Handback DateTime); zh-cn keeps the "not yet" placeholder.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$sheetsInfo = @(
    @{
        Name = "zh-cn"
        Row2Md = "e1d3f8d7-1db2-487e-b517-662f525dc137.md"
        Row2MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/f990e0bc33df199382549f31a369033e53c2d551/e2e/e1d3f8d7-1db2-487e-b517-662f525dc137.md"
        Row2Xlf = "e1d3f8d7-1db2-487e-b517-662f525dc137.5cb2d957b4a36af44e2b22bcb6d05d006b4b7d50.zh-cn.xlf"
        Row2XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b8e66a6dcd6a9c21ece7ccefb34363a83396feb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/e1d3f8d7-1db2-487e-b517-662f525dc137.5cb2d957b4a36af44e2b22bcb6d05d006b4b7d50.zh-cn.xlf"
        Row3Md = "e900ff5d-9406-4620-b0b0-4ecc073d7efd.md"
        Row3MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/f990e0bc33df199382549f31a369033e53c2d551/e2e/e900ff5d-9406-4620-b0b0-4ecc073d7efd.md"
        Row3Xlf = "e900ff5d-9406-4620-b0b0-4ecc073d7efd.1032a061caabb39fc32cf8f37d0bd5a03243da5d.zh-cn.xlf"
        Row3XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b8e66a6dcd6a9c21ece7ccefb34363a83396feb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/e900ff5d-9406-4620-b0b0-4ecc073d7efd.1032a061caabb39fc32cf8f37d0bd5a03243da5d.zh-cn.xlf"
        HandbackDateTime = $null
    },
    @{
        Name = "de-de"
        Row2Md = "e1d3f8d7-1db2-487e-b517-662f525dc137.md"
        Row2MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/f990e0bc33df199382549f31a369033e53c2d551/e2e/e1d3f8d7-1db2-487e-b517-662f525dc137.md"
        Row2Xlf = "e1d3f8d7-1db2-487e-b517-662f525dc137.5cb2d957b4a36af44e2b22bcb6d05d006b4b7d50.de-de.xlf"
        Row2XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/def380e91d845a9fab64f894f7b983dedc0805fc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/e1d3f8d7-1db2-487e-b517-662f525dc137.5cb2d957b4a36af44e2b22bcb6d05d006b4b7d50.de-de.xlf"
        Row3Md = "e900ff5d-9406-4620-b0b0-4ecc073d7efd.md"
        Row3MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/f990e0bc33df199382549f31a369033e53c2d551/e2e/e900ff5d-9406-4620-b0b0-4ecc073d7efd.md"
        Row3Xlf = "e900ff5d-9406-4620-b0b0-4ecc073d7efd.1032a061caabb39fc32cf8f37d0bd5a03243da5d.de-de.xlf"
        Row3XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/def380e91d845a9fab64f894f7b983dedc0805fc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/e900ff5d-9406-4620-b0b0-4ecc073d7efd.1032a061caabb39fc32cf8f37d0bd5a03243da5d.de-de.xlf"
        HandbackDateTime = "2016-01-17 07:34:43"
    }
)

# RGB(100,149,237) / #6495ED packed as BGR for the Font.Color OLE property,
# matching the workbook's existing custom "hyperlink" text color.
$linkColor = 15570276

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # --- Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Range("B2").Value = $newStatus
    $ws.Range("B3").Value = $newStatus

    # --- New "Latest Target File" (E) / "Latest Handback File" (F) links, row 2
    $ws.Hyperlinks.Add($ws.Range("E2"), $info.Row2MdUrl, "", "", $info.Row2Md)
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = $linkColor
    $ws.Range("E2").Font.Name = "Calibri"
    $ws.Range("E2").Font.Size = 11

    $ws.Hyperlinks.Add($ws.Range("F2"), $info.Row2XlfUrl, "", "", $info.Row2Xlf)
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = $linkColor
    $ws.Range("F2").Font.Name = "Calibri"
    $ws.Range("F2").Font.Size = 11

    # --- New "Latest Target File" (E) / "Latest Handback File" (F) links, row 3
    $ws.Hyperlinks.Add($ws.Range("E3"), $info.Row3MdUrl, "", "", $info.Row3Md)
    $ws.Range("E3").Font.Underline = $true
    $ws.Range("E3").Font.Color = $linkColor
    $ws.Range("E3").Font.Name = "Calibri"
    $ws.Range("E3").Font.Size = 11

    $ws.Hyperlinks.Add($ws.Range("F3"), $info.Row3XlfUrl, "", "", $info.Row3Xlf)
    $ws.Range("F3").Font.Underline = $true
    $ws.Range("F3").Font.Color = $linkColor
    $ws.Range("F3").Font.Name = "Calibri"
    $ws.Range("F3").Font.Size = 11

    # --- Latest Handback DateTime (G): only de-de got an actual handback stamp
    if ($info.HandbackDateTime) {
        $ws.Range("G2").Value = $info.HandbackDateTime
        $ws.Range("G3").Value = $info.HandbackDateTime
    }
}
